# Generate Report for handoff
#
# The source file "0cf96b12-b1ad-4b51-a6f3-3e6f07037073.md" was re-run through
# the handoff pipeline under a new id ("efd2c8ee-ff4c-4a8c-8790-f7dfe28ea8c8.md")
# and a second source file ("ffffdccbaa40-3d5d-471c-9c12-bb1fa08293e5.md") was
# added to the handoff batch. Both are now "Ready for handoff" / "Include"
# with real handoff file + datetime info recorded, and the
# ".localization-config" bookkeeping row is pushed down to make room.

$wb = $excel.ActiveWorkbook

$zhHandoffFile = "efd2c8ee-ff4c-4a8c-8790-f7dfe28ea8c8.70ae79f16e234a9eaa7e49d9fb4eb7ffc93e65e3.zh-cn.xlf"
$deHandoffFile = "efd2c8ee-ff4c-4a8c-8790-f7dfe28ea8c8.70ae79f16e234a9eaa7e49d9fb4eb7ffc93e65e3.de-de.xlf"

$mdUrlBase = "https://github.com/OpenLocalizationTest/oltest/blob/a9e99c6bbb63e0de4a742e20f3837f7b1dc9e6a9/e2e/"
$cfgUrl = "https://github.com/OpenLocalizationTest/oltest/blob/5611cd4322054c5c5d0d1a046a629b9394a11348/.localization-config"

# ---------------------------------------------------------------------------
# Overview sheet
# ---------------------------------------------------------------------------
$ws = $wb.Worksheets.Item("Overview")

$ws.Range("A1").Hyperlinks.Delete()
$ws.Rows.Item(3).Insert()

$ws.Range("A2").Value = "efd2c8ee-ff4c-4a8c-8790-f7dfe28ea8c8.md"
$ws.Range("B2").Value = "Ready for handoff"
$ws.Range("C2").Value = "Ready for handoff"

$ws.Range("A3").Value = "ffffdccbaa40-3d5d-471c-9c12-bb1fa08293e5.md"
$ws.Range("B3").Value = "Ready for handoff"
$ws.Range("C3").Value = "Ready for handoff"

$ws.Range("A4").Value = ".localization-config"
$ws.Range("B4").Value = "Not to be localized"
$ws.Range("C4").Value = "Not to be localized"

$ws.Hyperlinks.Add($ws.Range("A2"), "$mdUrlBase" + "efd2c8ee-ff4c-4a8c-8790-f7dfe28ea8c8.md", "", "", "efd2c8ee-ff4c-4a8c-8790-f7dfe28ea8c8.md")
$ws.Hyperlinks.Add($ws.Range("A3"), "$mdUrlBase" + "ffffdccbaa40-3d5d-471c-9c12-bb1fa08293e5.md", "", "", "ffffdccbaa40-3d5d-471c-9c12-bb1fa08293e5.md")
$ws.Hyperlinks.Add($ws.Range("A4"), $cfgUrl, "", "", ".localization-config")

# ---------------------------------------------------------------------------
# zh-cn sheet
# ---------------------------------------------------------------------------
$ws = $wb.Worksheets.Item("zh-cn")

$ws.Range("A1").Hyperlinks.Delete()
$ws.Rows.Item(3).Insert()

$ws.Range("A2").Value = "efd2c8ee-ff4c-4a8c-8790-f7dfe28ea8c8.md"
$ws.Range("B2").Value = "Ready for handoff"
$ws.Range("C2").Value = $zhHandoffFile
$ws.Range("D2").Value = "2016-01-21 02:52:58"
$ws.Range("G2").Value = "0001-01-01 00:00:00"
$ws.Range("H2").Value = "Include"

$ws.Range("A3").Value = "ffffdccbaa40-3d5d-471c-9c12-bb1fa08293e5.md"
$ws.Range("B3").Value = "Ready for handoff"
$ws.Range("C3").Value = $zhHandoffFile
$ws.Range("D3").Value = "2016-01-21 02:52:58"
$ws.Range("G3").Value = "0001-01-01 00:00:00"
$ws.Range("H3").Value = "Include"

$ws.Range("A4").Value = ".localization-config"
$ws.Range("B4").Value = "Not to be localized"
$ws.Range("D4").Value = "0001-01-01 00:00:00"
$ws.Range("G4").Value = "0001-01-01 00:00:00"
$ws.Range("H4").Value = "Ignored"

$ws.Hyperlinks.Add($ws.Range("A2"), "$mdUrlBase" + "efd2c8ee-ff4c-4a8c-8790-f7dfe28ea8c8.md", "", "", "efd2c8ee-ff4c-4a8c-8790-f7dfe28ea8c8.md")
$ws.Hyperlinks.Add($ws.Range("C2"), "$mdUrlBase" + "$zhHandoffFile", "", "", $zhHandoffFile)
$ws.Hyperlinks.Add($ws.Range("A3"), "$mdUrlBase" + "ffffdccbaa40-3d5d-471c-9c12-bb1fa08293e5.md", "", "", "ffffdccbaa40-3d5d-471c-9c12-bb1fa08293e5.md")
$ws.Hyperlinks.Add($ws.Range("C3"), "$mdUrlBase" + "$zhHandoffFile", "", "", $zhHandoffFile)
$ws.Hyperlinks.Add($ws.Range("A4"), $cfgUrl, "", "", ".localization-config")

# ---------------------------------------------------------------------------
# de-de sheet
# ---------------------------------------------------------------------------
$ws = $wb.Worksheets.Item("de-de")

$ws.Range("A1").Hyperlinks.Delete()
$ws.Rows.Item(3).Insert()

$ws.Range("A2").Value = "efd2c8ee-ff4c-4a8c-8790-f7dfe28ea8c8.md"
$ws.Range("B2").Value = "Ready for handoff"
$ws.Range("C2").Value = $deHandoffFile
$ws.Range("D2").Value = "2016-01-21 02:53:10"
$ws.Range("G2").Value = "0001-01-01 00:00:00"
$ws.Range("H2").Value = "Include"

$ws.Range("A3").Value = "ffffdccbaa40-3d5d-471c-9c12-bb1fa08293e5.md"
$ws.Range("B3").Value = "Ready for handoff"
$ws.Range("C3").Value = $deHandoffFile
$ws.Range("D3").Value = "2016-01-21 02:53:10"
$ws.Range("G3").Value = "0001-01-01 00:00:00"
$ws.Range("H3").Value = "Include"

$ws.Range("A4").Value = ".localization-config"
$ws.Range("B4").Value = "Not to be localized"
$ws.Range("D4").Value = "0001-01-01 00:00:00"
$ws.Range("G4").Value = "0001-01-01 00:00:00"
$ws.Range("H4").Value = "Ignored"

$ws.Hyperlinks.Add($ws.Range("A2"), "$mdUrlBase" + "efd2c8ee-ff4c-4a8c-8790-f7dfe28ea8c8.md", "", "", "efd2c8ee-ff4c-4a8c-8790-f7dfe28ea8c8.md")
$ws.Hyperlinks.Add($ws.Range("C2"), "$mdUrlBase" + "$deHandoffFile", "", "", $deHandoffFile)
$ws.Hyperlinks.Add($ws.Range("A3"), "$mdUrlBase" + "ffffdccbaa40-3d5d-471c-9c12-bb1fa08293e5.md", "", "", "ffffdccbaa40-3d5d-471c-9c12-bb1fa08293e5.md")
$ws.Hyperlinks.Add($ws.Range("C3"), "$mdUrlBase" + "$deHandoffFile", "", "", $deHandoffFile)
$ws.Hyperlinks.Add($ws.Range("A4"), $cfgUrl, "", "", ".localization-config")

Write-Output "Report generated for handoff"
